$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A35").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value = 46008
$ws.Range("B36").Value = 61

$ws.Range("A36:B36").Select()
